$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Insert two new worksheets after "missing_values" (4th sheet):
#   missing_axis_name
#   narrow_2d
# ---------------------------------------------------------------------
$wsMissingValues = $wb.Worksheets.Item(4)

$wsMissingAxis = $wb.Worksheets.Add($null, $wsMissingValues)
$wsMissingAxis.Name = "missing_axis_name"

$wsNarrow2d = $wb.Worksheets.Add($null, $wsMissingAxis)
$wsNarrow2d.Name = "narrow_2d"

# ---------------------------------------------------------------------
# Fill "missing_axis_name" sheet
#   a     b   c0  c1
#   a0    b0  0   1
#   a0    b1  2   3
#   a1    b0  4   5
#   a1    b1  6   7
# ---------------------------------------------------------------------
$wsMissingAxis.Range("A1").Value = "a"
$wsMissingAxis.Range("B1").Value = "b"
$wsMissingAxis.Range("C1").Value = "c0"
$wsMissingAxis.Range("D1").Value = "c1"

$wsMissingAxis.Range("A2").Value = "a0"
$wsMissingAxis.Range("B2").Value = "b0"
$wsMissingAxis.Range("C2").Value = 0
$wsMissingAxis.Range("D2").Value = 1

$wsMissingAxis.Range("A3").Value = "a0"
$wsMissingAxis.Range("B3").Value = "b1"
$wsMissingAxis.Range("C3").Value = 2
$wsMissingAxis.Range("D3").Value = 3

$wsMissingAxis.Range("A4").Value = "a1"
$wsMissingAxis.Range("B4").Value = "b0"
$wsMissingAxis.Range("C4").Value = 4
$wsMissingAxis.Range("D4").Value = 5

$wsMissingAxis.Range("A5").Value = "a1"
$wsMissingAxis.Range("B5").Value = "b1"
$wsMissingAxis.Range("C5").Value = 6
$wsMissingAxis.Range("D5").Value = 7

# Column A gets a distinct font + vertical-center alignment. Build the
# format on A1 only, then clone it onto A2:A5 via a formats-only paste so
# we don't leave behind unused intermediate style records.
$wsMissingAxis.Range("A1").Font.Name = "Calibri"
$wsMissingAxis.Range("A1").VerticalAlignment = -4108
$wsMissingAxis.Range("A1").Copy()
$wsMissingAxis.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsMissingAxis.Range("G18").Select()

# ---------------------------------------------------------------------
# Fill "narrow_2d" sheet (long/narrow format: a, b, value)
#   a  b   value
#   1  b0  0
#   1  b1  1
#   2  b0  2
#   2  b1  3
#   3  b0  4
#   3  b1  5
# ---------------------------------------------------------------------
$wsNarrow2d.Range("A1").Value = "a"
$wsNarrow2d.Range("B1").Value = "b"
$wsNarrow2d.Range("C1").Value = "value"

$wsNarrow2d.Range("A2").Value = 1
$wsNarrow2d.Range("B2").Value = "b0"
$wsNarrow2d.Range("C2").Value = 0

$wsNarrow2d.Range("A3").Value = 1
$wsNarrow2d.Range("B3").Value = "b1"
$wsNarrow2d.Range("C3").Value = 1

$wsNarrow2d.Range("A4").Value = 2
$wsNarrow2d.Range("B4").Value = "b0"
$wsNarrow2d.Range("C4").Value = 2

$wsNarrow2d.Range("A5").Value = 2
$wsNarrow2d.Range("B5").Value = "b1"
$wsNarrow2d.Range("C5").Value = 3

$wsNarrow2d.Range("A6").Value = 3
$wsNarrow2d.Range("B6").Value = "b0"
$wsNarrow2d.Range("C6").Value = 4

$wsNarrow2d.Range("A7").Value = 3
$wsNarrow2d.Range("B7").Value = "b1"
$wsNarrow2d.Range("C7").Value = 5

# narrow_2d ends up being the active sheet/tab
$wsNarrow2d.Range("F17").Select()

# ---------------------------------------------------------------------
# The previously-active last sheet ("int_labels") is no longer selected;
# selecting a cell elsewhere above already moved the "active sheet" flag
# off of it.
# ---------------------------------------------------------------------
